$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (merged A1:D1) from 2024-04-24 to 2024-05-24
$ws.Range("A1").Value = 45436

# Update prices in D28 and D29
$ws.Range("D28").Value = 47157
$ws.Range("D29").Value = 31460
